$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExpenseRequest")

# Update the error message text (remove trailing period)
$ws.Range("L2").Value = "Complete this field"

# Widen column L to fit the new text
$ws.Columns.Item(12).ColumnWidth = 22.8

# Update the active cell selection on the sheet
$ws.Range("L15").Select()
